# Generate Report for Handback
# The ad1f4133-... file has been handed back (in sync with en-US),
# so update its Status from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet, and stamp the
# Latest Handback DateTime for the per-locale sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Overview sheet: row for ad1f4133-...md is row 3; zh-cn (B) and de-de (C) status columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# zh-cn sheet: row 3 is the ad1f4133-...md entry. Column C = Status,
# column H = Latest Handback DateTime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("H3").Value = "2016-03-19 06:37:00"

# de-de sheet: row 3 is the ad1f4133-...md entry. Column C = Status,
# column H = Latest Handback DateTime.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("H3").Value = "2016-03-19 06:37:06"
